# Refresh of the "cryptos" price list: updates Price (D) / Volume(1h) (E)
# figures for most rows, and for rows 17/18 swaps BitcoinCash <-> WrappedBTC
# (coin name, link, price, volume) since their rank order changed.
#
# Column D values are entered with a leading apostrophe so Excel stores them
# as literal text (matching the source data, which uses a dotted thousands
# style like "26.744.98" and keeps trailing zeros like "1.00") instead of
# auto-converting them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.744.98"
$ws.Range("E2").Value = "  +2.00%  "
$ws.Range("D3").Value = "'1.624.41"
$ws.Range("E3").Value = "  +2.31%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'214.56"
$ws.Range("E5").Value = "  +1.17%  "
$ws.Range("E6").Value = "  +0.96%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("E8").Value = "  +0.47%  "
$ws.Range("E9").Value = "  +0.86%  "
$ws.Range("D10").Value = "'19.42"
$ws.Range("E10").Value = "  +0.50%  "
$ws.Range("E11").Value = "  +1.30%  "
$ws.Range("D12").Value = "'1.853.23"
$ws.Range("E12").Value = "  +2.34%  "
$ws.Range("D13").Value = "'1.629.66"
$ws.Range("E13").Value = "  +2.66%  "
$ws.Range("E14").Value = "  +1.38%  "
$ws.Range("D15").Value = "'65.08"
$ws.Range("E15").Value = "  +1.34%  "
$ws.Range("D16").Value = "'0.514"
$ws.Range("E16").Value = "  -1.12%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "'26.744.53"
$ws.Range("E17").Value = "  +1.98%  "
$ws.Range("B18").Value = "BitcoinCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D18").Value = "'235.02"
$ws.Range("E18").Value = "  +10.31%  "
$ws.Range("D19").Value = "'7.76"
$ws.Range("E19").Value = "  +5.15%  "
$ws.Range("E20").Value = "  +0.41%  "
$ws.Range("D21").Value = "'0.999"
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("D22").Value = "'4.41"
$ws.Range("E22").Value = "  +3.42%  "
$ws.Range("E23").Value = "  +3.76%  "
$ws.Range("D24").Value = "'9.12"
$ws.Range("E24").Value = "  +1.50%  "
$ws.Range("D25").Value = "'145.93"
$ws.Range("E25").Value = "  +1.69%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("E27").Value = "  +0.88%  "
$ws.Range("E28").Value = "  +2.60%  "
$ws.Range("D29").Value = "'15.72"
$ws.Range("E29").Value = "  +3.64%  "
$ws.Range("E30").Value = "  +0.55%  "
$ws.Range("E31").Value = "  +1.07%  "
$ws.Range("E32").Value = "  +2.11%  "
$ws.Range("D33").Value = "'1.473.46"
$ws.Range("E33").Value = "  +10.07%  "
$ws.Range("E34").Value = "  +2.33%  "
$ws.Range("E36").Value = "  +2.06%  "
$ws.Range("E38").Value = "  +0.60%  "
$ws.Range("E39").Value = "  +2.19%  "
$ws.Range("D40").Value = "'5.96"
$ws.Range("E40").Value = "  +3.65%  "
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("E42").Value = "  +3.13%  "
$ws.Range("D43").Value = "'0.951"
$ws.Range("E43").Value = "  -2.09%  "
$ws.Range("D44").Value = "'1.764.38"
$ws.Range("E44").Value = "  +2.42%  "
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("D46").Value = "'62.12"
$ws.Range("D47").Value = "'88.51"
$ws.Range("E47").Value = "  +3.27%  "
$ws.Range("E48").Value = "  +2.32%  "
$ws.Range("D49").Value = "'0.0505"
$ws.Range("E49").Value = "  +0.66%  "
$ws.Range("D50").Value = "'0.0967"
$ws.Range("E50").Value = "  -0.97%  "
$ws.Range("D51").Value = "'7.50"
$ws.Range("E51").Value = "  +1.78%  "
